# Update countries & provincias Spain
# Applies the 30-Apr-2020 20:22 data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-rank "Irak": its case count overtook Croacia / Ghana / Armenia /
# Uzbekistan, so those four rows shift down one slot and Irak takes row 68.
# (Row 73 / Camerun is unaffected.)
$ws.Range("A68").Value = "Irak"
$ws.Range("A69").Value = "Croacia"
$ws.Range("A70").Value = "Ghana"
$ws.Range("A71").Value = "Armenia"
$ws.Range("A72").Value = "Uzbekistan"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1078476
$ws.Range("C4").Value = 14282
$ws.Range("E4").Value = 866255
$ws.Range("G4").Value = 880
$ws.Range("H4").Value = 62535

# --- Row 15 ---
$ws.Range("B15").Value = 53013
$ws.Range("C15").Value = 1416
$ws.Range("D15").Value = 21187
$ws.Range("E15").Value = 28646
$ws.Range("G15").Value = 184
$ws.Range("H15").Value = 3180

# --- Row 37 ---
$ws.Range("E37").Value = 7512
$ws.Range("G37").Value = 18
$ws.Range("H37").Value = 711

# --- Row 44 ---
$ws.Range("E44").Value = 7496
$ws.Range("F44").Value = 37
$ws.Range("G44").Value = 3
$ws.Range("H44").Value = 210

# --- Row 68: now Irak ---
$ws.Range("B68").Value = 2085
$ws.Range("C68").Value = 82
$ws.Range("D68").Value = 1375
$ws.Range("E68").Value = 617
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 93

# --- Row 69: now Croacia ---
$ws.Range("B69").Value = 2076
$ws.Range("C69").Value = 14
$ws.Range("D69").Value = 1348
$ws.Range("E69").Value = 659
$ws.Range("F69").Value = 20
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 69

# --- Row 70: now Ghana ---
$ws.Range("B70").Value = 2074
$ws.Range("C70").Value = 403
$ws.Range("D70").Value = 212
$ws.Range("E70").Value = 1845
$ws.Range("F70").Value = 4
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 17

# --- Row 71: now Armenia ---
$ws.Range("B71").Value = 2066
$ws.Range("C71").Value = 134
$ws.Range("D71").Value = 929
$ws.Range("E71").Value = 1105
$ws.Range("F71").Value = 10
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = 32

# --- Row 72: now Uzbekistan ---
$ws.Range("B72").Value = 2017
$ws.Range("C72").Value = 15
$ws.Range("D72").Value = 1133
$ws.Range("E72").Value = 875
$ws.Range("F72").Value = 8
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 9

# --- Timestamp footer ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 20:22"
